$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-39: update the date serial value from 45190 to 45192
$ws.Range("C2:C39").Value = 45192
